$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.231415721352562
$ws.Cells.Item(2, 3).Value = 0.31266005797562
$ws.Cells.Item(2, 4).Value = 0.0792997928117245
$ws.Cells.Item(2, 5).Value = 0.4210711681848807
$ws.Cells.Item(2, 7).Value = 0.002393885002521618
$ws.Cells.Item(2, 15).Value = 2.150681531974698
$ws.Cells.Item(3, 2).Value = 1.086345546104383
$ws.Cells.Item(3, 3).Value = 0.2732052880114111
$ws.Cells.Item(3, 4).Value = 0.07181007949643003
$ws.Cells.Item(3, 5).Value = 0.3671310163277752
$ws.Cells.Item(3, 7).Value = 0.002397947023809593
$ws.Cells.Item(3, 15).Value = 2.131775290707424
$ws.Cells.Item(4, 2).Value = 0.9972292380033423
$ws.Cells.Item(4, 3).Value = 0.248887350652609
$ws.Cells.Item(4, 4).Value = 0.06725008652669828
$ws.Cells.Item(4, 5).Value = 0.3341271889903936
$ws.Cells.Item(4, 7).Value = 0.002400568940967025
$ws.Cells.Item(4, 15).Value = 2.122666070498553
$ws.Cells.Item(5, 2).Value = 0.9609034933564544
$ws.Cells.Item(5, 3).Value = 0.2389545559492774
$ws.Cells.Item(5, 4).Value = 0.06540147747342928
$ws.Cells.Item(5, 5).Value = 0.32070449458773
$ws.Cells.Item(5, 7).Value = 0.002401669647512616
$ws.Cells.Item(5, 15).Value = 2.119577385641009
$ws.Cells.Item(6, 2).Value = 0.9548710381549768
$ws.Cells.Item(6, 3).Value = 0.2373038387319752
$ws.Cells.Item(6, 4).Value = 0.06509509545070102
$ws.Cells.Item(6, 5).Value = 0.3184772039549415
$ws.Cells.Item(6, 7).Value = 0.002401854370366179
$ws.Cells.Item(6, 15).Value = 2.119102015425625
$ws.Cells.Item(7, 2).Value = 0.9967393764796952
$ws.Cells.Item(7, 3).Value = 0.2487534864617089
$ws.Cells.Item(7, 4).Value = 0.06722511671141262
$ws.Cells.Item(7, 5).Value = 0.3339460611963432
$ws.Cells.Item(7, 7).Value = 0.002400583654784756
$ws.Cells.Item(7, 15).Value = 2.122621898043462
$ws.Cells.Item(8, 2).Value = 1.181404672473377
$ws.Cells.Item(8, 3).Value = 0.2990753867852902
$ws.Cells.Item(8, 4).Value = 0.07670923311877687
$ws.Cells.Item(8, 5).Value = 0.4024471079228107
$ws.Cells.Item(8, 7).Value = 0.002395259126405013
$ws.Cells.Item(8, 15).Value = 2.143640864449338
$ws.Cells.Item(9, 2).Value = 1.543188713654047
$ws.Cells.Item(9, 3).Value = 0.3970172172007551
$ws.Cells.Item(9, 4).Value = 0.09562085764689243
$ws.Cells.Item(9, 5).Value = 0.5378121374655933
$ws.Cells.Item(9, 7).Value = 0.002385826775525977
$ws.Cells.Item(9, 15).Value = 2.204924301221354
$ws.Cells.Item(10, 2).Value = 1.808804166832942
$ws.Cells.Item(10, 3).Value = 0.4685257348454002
$ws.Cells.Item(10, 4).Value = 0.1097161047633506
$ws.Cells.Item(10, 5).Value = 0.6380675335471295
$ws.Cells.Item(10, 7).Value = 0.002379504705924882
$ws.Cells.Item(10, 15).Value = 2.262514896549391
$ws.Cells.Item(11, 2).Value = 1.929606308001439
$ws.Cells.Item(11, 3).Value = 0.5009602504061945
$ws.Cells.Item(11, 4).Value = 0.1161742561311456
$ws.Cells.Item(11, 5).Value = 0.6838887419247186
$ws.Cells.Item(11, 7).Value = 0.002376759074711203
$ws.Cells.Item(11, 15).Value = 2.291515586040532
$ws.Cells.Item(12, 2).Value = 1.975347023833137
$ws.Cells.Item(12, 3).Value = 0.5132286276538025
$ws.Cells.Item(12, 4).Value = 0.1186265792390202
$ws.Cells.Item(12, 5).Value = 0.7012738501878033
$ws.Cells.Item(12, 7).Value = 0.002375737995295387
$ws.Cells.Item(12, 15).Value = 2.302906177714021
$ws.Cells.Item(13, 2).Value = 1.965496146150826
$ws.Cells.Item(13, 3).Value = 0.5105870329843469
$ws.Cells.Item(13, 4).Value = 0.1180981248624136
$ws.Cells.Item(13, 5).Value = 0.6975281167607363
$ws.Cells.Item(13, 7).Value = 0.002375957076098406
$ws.Cells.Item(13, 15).Value = 2.300434743965411
$ws.Cells.Item(14, 2).Value = 1.933369515497361
$ws.Cells.Item(14, 3).Value = 0.5019698562451822
$ws.Cells.Item(14, 4).Value = 0.1163758739998713
$ws.Cells.Item(14, 5).Value = 0.6853183356484607
$ws.Cells.Item(14, 7).Value = 0.002376674697048067
$ws.Cells.Item(14, 15).Value = 2.292444474084789
$ws.Cells.Item(15, 2).Value = 1.913690436597506
$ws.Cells.Item(15, 3).Value = 0.4966897747545431
$ws.Cells.Item(15, 4).Value = 0.11532182987753
$ws.Cells.Item(15, 5).Value = 0.6778439520548005
$ws.Cells.Item(15, 7).Value = 0.002377116683901737
$ws.Cells.Item(15, 15).Value = 2.287603591694221
$ws.Cells.Item(16, 2).Value = 1.800908866074451
$ws.Cells.Item(16, 3).Value = 0.4664041325926291
$ws.Cells.Item(16, 4).Value = 0.1092949887043915
$ws.Cells.Item(16, 5).Value = 0.6350775618046498
$ws.Cells.Item(16, 7).Value = 0.002379686752536204
$ws.Cells.Item(16, 15).Value = 2.260676509533965
$ws.Cells.Item(17, 2).Value = 1.731713600478486
$ws.Cells.Item(17, 3).Value = 0.4478003764239133
$ws.Cells.Item(17, 4).Value = 0.1056096257586603
$ws.Cells.Item(17, 5).Value = 0.6088986639848741
$ws.Cells.Item(17, 7).Value = 0.00238129670563548
$ws.Cells.Item(17, 15).Value = 2.244879412346421
$ws.Cells.Item(18, 2).Value = 1.691911734284304
$ws.Cells.Item(18, 3).Value = 0.4370910506516452
$ws.Cells.Item(18, 4).Value = 0.1034942390121074
$ws.Cells.Item(18, 5).Value = 0.5938612787813611
$ws.Cells.Item(18, 7).Value = 0.002382234979788499
$ws.Cells.Item(18, 15).Value = 2.236056777236826
$ws.Cells.Item(19, 2).Value = 1.678435077972438
$ws.Cells.Item(19, 3).Value = 0.4334635282927479
$ws.Cells.Item(19, 4).Value = 0.1027787464475693
$ws.Cells.Item(19, 5).Value = 0.5887732342775109
$ws.Cells.Item(19, 7).Value = 0.00238255477438737
$ws.Cells.Item(19, 15).Value = 2.233114657968969
$ws.Cells.Item(20, 2).Value = 1.739079825048179
$ws.Cells.Item(20, 3).Value = 0.4497817044778003
$ws.Cells.Item(20, 4).Value = 0.1060014892192953
$ws.Cells.Item(20, 5).Value = 0.6116833599679126
$ws.Cells.Item(20, 7).Value = 0.002381124054121758
$ws.Cells.Item(20, 15).Value = 2.246533734657731
$ws.Cells.Item(21, 2).Value = 1.942806009375033
$ws.Cells.Item(21, 3).Value = 0.5045013078747047
$ws.Cells.Item(21, 4).Value = 0.1168815563230368
$ws.Cells.Item(21, 5).Value = 0.6889037099066115
$ws.Cells.Item(21, 7).Value = 0.002376463409612306
$ws.Cells.Item(21, 15).Value = 2.294780273712036
$ws.Cells.Item(22, 2).Value = 2.075926832446669
$ws.Cells.Item(22, 3).Value = 0.5401828059275431
$ws.Cells.Item(22, 4).Value = 0.1240317786553504
$ws.Cells.Item(22, 5).Value = 0.7395688968093879
$ws.Cells.Item(22, 7).Value = 0.002373525958637803
$ws.Cells.Item(22, 15).Value = 2.328696820530354
$ws.Cells.Item(23, 2).Value = 2.004880261455639
$ws.Cells.Item(23, 3).Value = 0.5211463855555394
$ws.Cells.Item(23, 4).Value = 0.1202119172783114
$ws.Cells.Item(23, 5).Value = 0.712508962967803
$ws.Cells.Item(23, 7).Value = 0.002375083834333189
$ws.Cells.Item(23, 15).Value = 2.31037484118508
$ws.Cells.Item(24, 2).Value = 1.735749619014484
$ws.Cells.Item(24, 3).Value = 0.4488859889722789
$ws.Cells.Item(24, 4).Value = 0.1058243172401205
$ws.Cells.Item(24, 5).Value = 0.6104243579993494
$ws.Cells.Item(24, 7).Value = 0.002381202070183441
$ws.Cells.Item(24, 15).Value = 2.245785008185265
$ws.Cells.Item(25, 2).Value = 1.445351244785058
$ws.Cells.Item(25, 3).Value = 0.3706004182140532
$ws.Cells.Item(25, 4).Value = 0.09047020645515147
$ws.Cells.Item(25, 5).Value = 0.5010639066571798
$ws.Cells.Item(25, 7).Value = 0.002388271203781306
$ws.Cells.Item(25, 15).Value = 2.186163050303662
